# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
#
# This appends the new week's per-play logs to the four running-log cells on
# the "YDS" sheet and the six running-log cells on the "ST" sheet, then
# updates the season-total cells on "OFF", "DEF", "ST", "TURNS" and "PEN"
# to reflect the new week's contribution.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append Week 17 per-play yardage logs (space separated ints)
# ---------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value2 + " 3 -3 4 6 4 18 24 3 0 1 1 14 2 0 4 9 -1 2 8 5 7 1 4 3 1"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value2 + " 6 2 6 4 6 9 2 9 0 0 2 -4 1 -1 3 5 4 4 8 3 2 1 1 2 5 2"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value2 + " 8 5 17 7 6 8 4 12 7 45 5 3 12 6 0 3"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value2 + " 14 5 7 3 11 -2 11 4 26 12 5 4 5 19 13 7 42 2 5 9 4 14 24"

# ---------------------------------------------------------------------
# OFF sheet: updated season totals (Home row 2, Road row 3)
# ---------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value = 204
$wsOFF.Range("D2").Value = 16
$wsOFF.Range("F2").Value = 83
$wsOFF.Range("G2").Value = 60
$wsOFF.Range("J2").Value = 36
$wsOFF.Range("L2").Value = 292
$wsOFF.Range("M2").Value = 176
$wsOFF.Range("O2").Value = 23
$wsOFF.Range("P2").Value = 15
$wsOFF.Range("Q2").Value = 550

$wsOFF.Range("C3").Value = 156
$wsOFF.Range("D3").Value = 10
$wsOFF.Range("E3").Value = 37
$wsOFF.Range("F3").Value = 77
$wsOFF.Range("G3").Value = 21
$wsOFF.Range("H3").Value = 24
$wsOFF.Range("I3").Value = 58
$wsOFF.Range("J3").Value = 45
$wsOFF.Range("N3").Value = 16

# ---------------------------------------------------------------------
# DEF sheet: updated season totals (Home row 2, Road row 3)
# ---------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value = 182
$wsDEF.Range("D2").Value = 9
$wsDEF.Range("E2").Value = 8
$wsDEF.Range("F2").Value = 61
$wsDEF.Range("G2").Value = 54
$wsDEF.Range("J2").Value = 27
$wsDEF.Range("L2").Value = 329
$wsDEF.Range("M2").Value = 215
$wsDEF.Range("O2").Value = 20
$wsDEF.Range("P2").Value = 11
$wsDEF.Range("Q2").Value = 586

$wsDEF.Range("B3").Value = 11
$wsDEF.Range("C3").Value = 196
$wsDEF.Range("E3").Value = 24
$wsDEF.Range("F3").Value = 117
$wsDEF.Range("G3").Value = 33
$wsDEF.Range("H3").Value = 23
$wsDEF.Range("I3").Value = 66
$wsDEF.Range("J3").Value = 57
$wsDEF.Range("N3").Value = 19

# ---------------------------------------------------------------------
# ST sheet: updated season totals + appended Week 17 logs
# ---------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 91
$wsST.Range("D2").Value = 56
$wsST.Range("F2").Value = 131
$wsST.Range("G2").Value = 126
$wsST.Range("L2").Value = 35
$wsST.Range("M2").Value = 26

$wsST.Range("B3").Value = 45

$wsST.Range("D3").Value = $wsST.Range("D3").Value2 + " 45 40 45 38 55"
$wsST.Range("B4").Value = $wsST.Range("B4").Value2 + " 64 63"
$wsST.Range("D4").Value = $wsST.Range("D4").Value2 + " 7 41 0 0 16"
$wsST.Range("B5").Value = $wsST.Range("B5").Value2 + " 34 9"
$wsST.Range("D5").Value = $wsST.Range("D5").Value2 + " 0 0"
$wsST.Range("B6").Value = $wsST.Range("B6").Value2 + " 20 0"

# ---------------------------------------------------------------------
# TURNS sheet: updated season totals
# ---------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("C2").Value = 11
$wsTURNS.Range("D2").Value = 7
$wsTURNS.Range("E2").Value = 9

$wsTURNS.Range("D3").Value = 6

# ---------------------------------------------------------------------
# PEN sheet: updated season totals
# ---------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B3").Value = 21

Write-Host "Week 17 data logged."
